# This commit ("G2-1774 Update apache POI") only reflects an Apache POI
# library version bump that was used to originally produce/re-save this
# fixture .pptx. Diffing the canonical OOXML shows every changed line is
# a pure XML attribute re-ordering (e.g. `xmlns:a/xmlns:p/xmlns:r` vs
# `xmlns:a/xmlns:r/xmlns:p`, `kern`/`sz` swapped, `idx`/`type` swapped,
# etc.) across presentation.xml, the slide layouts, the slide master and
# the theme - there is no change to any attribute *value*, no element
# added or removed, and no slide content touched at all
# (ppt/slides/slide1.xml does not even appear in the diff).
#
# That kind of byte-level re-serialization is an artifact of the XML
# writer inside the library that saved the file, not something that is
# reachable (or meaningful) through the PowerPoint object model - there
# is no COM property for "attribute order". So the faithful application
# of this change through PowerPoint automation is to touch nothing: no
# slide, shape, text, layout, master or theme property changes.

$p = $ppt.ActivePresentation
# No-op: confirm the presentation handle resolves; intentionally make no
# content changes, matching the semantically-empty diff above.
$null = $p.Slides.Count
